# Regenerate save_data column G ("K") values for rows 2-71.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$kValues = @{
    2 = 0
    3 = 2
    4 = 0
    5 = 0
    6 = 0
    7 = 1
    8 = 1
    9 = 1
    10 = 2
    11 = 1
    12 = 3
    13 = 1
    14 = 1
    15 = 1
    16 = 0
    17 = 1
    18 = 0
    19 = 0
    20 = 2
    21 = 2
    22 = 1
    23 = 0
    24 = 0
    25 = 0
    26 = 4
    27 = 1
    28 = 3
    29 = 0
    30 = 0
    31 = 1
    32 = 2
    33 = 0
    34 = 3
    35 = 2
    36 = 0
    37 = 1
    38 = 1
    39 = 1
    40 = 2
    41 = 0
    42 = 2
    43 = 0
    44 = 1
    45 = 1
    46 = 2
    47 = 2
    48 = 0
    49 = 1
    50 = 1
    51 = 0
    52 = 0
    53 = 0
    54 = 1
    55 = 1
    56 = 0
    57 = 2
    58 = 2
    59 = 0
    60 = 1
    61 = 1
    62 = 2
    63 = 0
    64 = 1
    65 = 2
    66 = 1
    67 = 0
    68 = 1
    69 = 1
    70 = 1
    71 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
